$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new group header "M_PL" merged across R1:Y1 ---
$ws.Range("R1:Y1").Merge() | Out-Null
$ws.Range("R1").Value = "M_PL"

# --- Row 2: new sub-headers for the M_PL group (R2:Y2), mirrors the pattern
#     already used for the M_%cit (B:I) and M_ETR (J:Q) groups ---
$ws.Range("R2").Value = "GFA - Sales"
$ws.Range("S2").Value = "GFA - Sales + Emp"
$ws.Range("T2").Value = "IMF - Sales"
$ws.Range("U2").Value = "IMF - Sales + Emp"
$ws.Range("V2").Value = "OECD (20%) - Sales"
$ws.Range("W2").Value = "OECD (20%) - Sales + Emp"
$ws.Range("X2").Value = "OECD - Sales"
$ws.Range("Y2").Value = "OECD - Sales + Emp"

# --- Apply the same formatting used by the rest of the header block
#     (bold font, centered/top aligned, thin border all around) ---
$rng = $ws.Range("R1:Y2")
$rng.Font.Bold = $true
$rng.HorizontalAlignment = -4108  # xlCenter
$rng.VerticalAlignment = -4160    # xlTop
$rng.Borders.LineStyle = 1

# --- Data rows 4-10: new profit figures in columns R:Y ---
$ws.Range("R4").Value = 40382619918
$ws.Range("S4").Value = 37273404592
$ws.Range("T4").Value = 30235087160
$ws.Range("U4").Value = 37431807140
$ws.Range("V4").Value = 66677809103
$ws.Range("W4").Value = 66677809103
$ws.Range("X4").Value = 66677809103
$ws.Range("Y4").Value = 66677809103

$ws.Range("R5").Value = 68341277913
$ws.Range("S5").Value = 68607173698
$ws.Range("T5").Value = 66666309306
$ws.Range("U5").Value = 66917200922
$ws.Range("V5").Value = 75243181211
$ws.Range("W5").Value = 75243181211
$ws.Range("X5").Value = 75243181211
$ws.Range("Y5").Value = 75243181211

$ws.Range("R6").Value = 16282609381
$ws.Range("S6").Value = 16245802182
$ws.Range("T6").Value = 15485301894
$ws.Range("U6").Value = 16245802182
$ws.Range("V6").Value = 21419969152
$ws.Range("W6").Value = 21419969152
$ws.Range("X6").Value = 21419969152
$ws.Range("Y6").Value = 21419969152

$ws.Range("R7").Value = 1403512928
$ws.Range("S7").Value = 1403512928
$ws.Range("T7").Value = 1403512928
$ws.Range("U7").Value = 1403512928
$ws.Range("V7").Value = 6471301295
$ws.Range("W7").Value = 6471301295
$ws.Range("X7").Value = 6471301295
$ws.Range("Y7").Value = 6471301295

$ws.Range("R8").Value = 885447038872
$ws.Range("S8").Value = 885447038872
$ws.Range("T8").Value = 885447038872
$ws.Range("U8").Value = 885447038872
$ws.Range("V8").Value = 885447038872
$ws.Range("W8").Value = 885447038872
$ws.Range("X8").Value = 885447038872
$ws.Range("Y8").Value = 885447038872

$ws.Range("S9").Value = 12096909667
$ws.Range("U9").Value = 12096909667
$ws.Range("V9").Value = 12956669707
$ws.Range("W9").Value = 12956669707
$ws.Range("X9").Value = 12956669707
$ws.Range("Y9").Value = 12956669707

$ws.Range("R10").Value = 1986738567
$ws.Range("S10").Value = 1195614519
$ws.Range("T10").Value = 1986738567
$ws.Range("U10").Value = 2032015673
$ws.Range("V10").Value = 2539331704
$ws.Range("W10").Value = 2539331704
$ws.Range("X10").Value = 2539331704
$ws.Range("Y10").Value = 2539331704
